$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data for rows 2-27 (columns A-S), reflecting the new block/trial numbering
# and refreshed stimulus assignments for the kitchens categorization task.
$data = @(
    @(11, 'categorization', 6, 2, 1, 271, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_syam3.png', 41.32432432432432, 26.2972972972973, 33.81081081081081, 37, 2, 2, 2),
    @(11, 'categorization', 6, 2, 2, 272, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_ra2nm.png', 70.75, 50.375, 60.5625, 40, 6, 6, 6),
    @(11, 'categorization', 6, 2, 3, 273, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_xr3up.png', 76.24444444444444, 55.88888888888889, 66.06666666666666, 45, 7, 7, 7),
    @(11, 'categorization', 6, 2, 4, 274, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_9oofc.png', 82.47619047619048, 65.5, 73.98809523809524, 42, 8, 8, 8),
    @(11, 'categorization', 6, 2, 5, 275, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_kq9s9.png', 62.30232558139535, 39.97674418604651, 51.13953488372093, 43, 4, 4, 4),
    @(11, 'categorization', 6, 2, 6, 276, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_swq34.png', 64.11363636363636, 43.04545454545455, 53.57954545454545, 44, 5, 5, 5),
    @(11, 'categorization', 6, 2, 7, 277, 'kitchens', 'bedrooms', 'target', $null, 'j', 'stimuli/img_5il0t.png', 48.09523809523809, 30.90476190476191, 39.5, 42, 2, 2, 2),
    @(11, 'categorization', 6, 2, 8, 278, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_24rt2.png', 55.26829268292683, 34.19512195121951, 44.73170731707317, 41, 3, 3, 3),
    @(11, 'categorization', 6, 2, 9, 279, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_vgh2g.png', 93.81395348837209, 78.27906976744185, 86.04651162790697, 43, 10, 10, 10),
    @(11, 'categorization', 6, 2, 10, 280, 'kitchens', 'bedrooms', 'distractor', $null, 'f', 'stimuli/img_twj5p.png', 67.71739130434783, 42.08695652173913, 54.90217391304348, 46, 4, 4, 4),
    @(11, 'categorization', 6, 2, 11, 281, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_2qhro.png', 81.73809523809524, 62.73809523809524, 72.23809523809524, 42, 8, 8, 8),
    @(11, 'categorization', 6, 2, 12, 282, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_iudc4.png', 73.625, 52.275, 62.95, 40, 6, 6, 6),
    @(11, 'categorization', 6, 2, 13, 283, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_rru0v.png', 56.45238095238095, 39.42857142857143, 47.94047619047619, 42, 4, 4, 4),
    @(11, 'categorization', 6, 2, 14, 284, 'kitchens', 'bedrooms', 'distractor', $null, 'f', 'stimuli/img_3h4c9.png', 85.47619047619048, 67.26190476190476, 76.36904761904762, 42, 9, 9, 9),
    @(11, 'categorization', 6, 2, 15, 285, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_j4ttn.png', 12.61904761904762, 11.42857142857143, 12.02380952380952, 42, 1, 1, 1),
    @(11, 'categorization', 6, 2, 16, 286, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_zxvl3.png', 68.78260869565217, 47.56521739130435, 58.17391304347827, 46, 5, 5, 5),
    @(11, 'categorization', 6, 2, 17, 287, 'kitchens', 'bedrooms', 'distractor', $null, 'f', 'stimuli/img_jp28n.png', 65.02564102564102, 44.97435897435897, 55, 39, 4, 4, 4),
    @(11, 'categorization', 6, 2, 18, 288, 'kitchens', 'bedrooms', 'target', $null, 'j', 'stimuli/img_d3t0o.png', 66.95121951219512, 42.92682926829269, 54.9390243902439, 41, 4, 4, 4),
    @(11, 'categorization', 6, 2, 19, 289, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_rg4in.png', 49.3695652173913, 30.21739130434782, 39.79347826086956, 46, 3, 3, 3),
    @(11, 'categorization', 6, 2, 20, 290, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_of8d6.png', 26.04878048780488, 19.14634146341463, 22.59756097560975, 41, 1, 1, 1),
    @(11, 'categorization', 6, 2, 21, 291, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_tn8ys.png', 86.70454545454545, 72.4090909090909, 79.55681818181819, 44, 10, 10, 10),
    @(11, 'categorization', 6, 2, 22, 292, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_amsgw.png', 86.08510638297872, 65.95744680851064, 76.02127659574468, 47, 9, 9, 9),
    @(11, 'categorization', 6, 2, 23, 293, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_vh7v8.png', 78.70454545454545, 59.63636363636363, 69.17045454545455, 44, 7, 7, 7),
    @(11, 'categorization', 6, 2, 24, 294, 'kitchens', 'bedrooms', 'target', $null, 'j', 'stimuli/img_91csq.png', 50.44736842105263, 28.34210526315789, 39.39473684210526, 38, 2, 2, 2),
    @(11, 'categorization', 6, 2, 25, 295, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_ac0ey.png', 86.62222222222222, 70.02222222222223, 78.32222222222222, 45, 9, 9, 9),
    @(11, 'categorization', 6, 2, 26, 296, 'kitchens', 'living_rooms', 'target', $null, 'j', 'stimuli/img_rych7.png', 30.4468085106383, 23.4468085106383, 26.9468085106383, 47, 2, 2, 2)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

# Remove the now-obsolete trailing row (previously row 28) and shrink dimension to A1:S27
$ws.Rows.Item(28).Delete()

Write-Host ("New used range: " + $ws.UsedRange.Address())
